# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.199.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.557.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.552.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.10%  "

# Row 8
$ws.Range("E8").Value = "  +0.97%  "

# Row 9
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("E10").Value = "  +9.44%  "

# Row 11
$ws.Range("E11").Value = "  +0.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000309"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.122.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.73%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.261.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.79%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.541.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.86%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.72%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "573.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.64%  "

# Row 21
$ws.Range("E21").Value = "  +0.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.76%  "

# Row 27
$ws.Range("E27").Value = "  -1.41%  "

# Row 28
$ws.Range("E28").Value = "  -2.87%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.00%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.43%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.96%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.43%  "

# Row 33
$ws.Range("E33").Value = "  +2.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.49%  "

# Row 35
$ws.Range("E35").Value = "  +21.58%  "

# Row 36
$ws.Range("E36").Value = "  +0.95%  "

# Row 37
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.404"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "521.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.77%  "

# Row 39
$ws.Range("E39").Value = "  +0.10%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.95%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.601.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.91%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0779"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.01%  "

# Row 43
$ws.Range("E43").Value = "  +3.70%  "

# Row 44
$ws.Range("E44").Value = "  +1.91%  "

# Row 45
$ws.Range("E45").Value = "  +3.74%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.85%  "

# Row 47
$ws.Range("E47").Value = "  -1.14%  "

# Row 48
$ws.Range("E48").Value = "  +3.24%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "

# Row 50
$ws.Range("E50").Value = "  +0.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.12%  "
